$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.769.18'
$ws.Range("E2").Value = '  +1.02%  '
$ws.Range("D3").Value = '1.825.71'
$ws.Range("E3").Value = '  +1.76%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.13'
$ws.Range("E5").Value = '  +1.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.576'
$ws.Range("E6").Value = '  +3.96%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '34.90'
$ws.Range("E8").Value = '  +7.88%  '
$ws.Range("E9").Value = '  +2.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0700'
$ws.Range("E10").Value = '  +1.08%  '
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("D12").Value = '2.089.29'
$ws.Range("E12").Value = '  +1.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.48'
$ws.Range("E13").Value = '  +4.09%  '
$ws.Range("D14").Value = '1.808.02'
$ws.Range("E14").Value = '  +0.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.649'
$ws.Range("E15").Value = '  +2.85%  '
$ws.Range("D16").Value = '34.755.26'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.35'
$ws.Range("E17").Value = '  +3.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.18'
$ws.Range("E18").Value = '  +1.29%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '248.10'
$ws.Range("E19").Value = '  +0.71%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0804'
$ws.Range("E20").Value = '  +0.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.65'
$ws.Range("E21").Value = '  +5.56%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.21'
$ws.Range("E23").Value = '  +1.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '174.22'
$ws.Range("E24").Value = '  +6.98%  '
$ws.Range("E25").Value = '  +1.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.50'
$ws.Range("E26").Value = '  +3.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.89'
$ws.Range("E27").Value = '  +3.06%  '
$ws.Range("E28").Value = '  +2.46%  '
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.02'
$ws.Range("E30").Value = '  +3.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0534'
$ws.Range("E31").Value = '  +2.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.87'
$ws.Range("E32").Value = '  +2.63%  '
$ws.Range("E33").Value = '  +1.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.87'
$ws.Range("E34").Value = '  +2.21%  '
$ws.Range("E35").Value = '  +0.57%  '
$ws.Range("D36").Value = '1.420.46'
$ws.Range("E36").Value = '  -1.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.687'
$ws.Range("E37").Value = '  +2.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.08'
$ws.Range("E38").Value = '  +2.12%  '
$ws.Range("E39").Value = '  +0.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '85.55'
$ws.Range("E40").Value = '  +1.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.88'
$ws.Range("E41").Value = '  +4.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.962'
$ws.Range("E42").Value = '  +3.03%  '
$ws.Range("E43").Value = '  +0.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.87'
$ws.Range("E44").Value = '  +0.53%  '
$ws.Range("E45").Value = '  +2.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0518'
$ws.Range("E46").Value = '  -1.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.12'
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("D48").Value = '1.989.24'
$ws.Range("E48").Value = '  +1.91%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.91'
$ws.Range("E49").Value = '  +0.22%  '
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  -0.11%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0128'
$ws.Range("E51").Value = '  -1.93%  '
